$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 79, pushing the existing data (rows 79-127) down
# by one row (they become rows 80-128).
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new price record.
$ws.Cells.Item(79, 1).Value = 10
$ws.Cells.Item(79, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(79, 3).Value = "La Araucanía"
$ws.Cells.Item(79, 4).Value = 45233
$ws.Cells.Item(79, 5).Value = 9
$ws.Cells.Item(79, 6).Value = 100112022
$ws.Cells.Item(79, 7).Value = "Arveja Verde"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 110
$ws.Cells.Item(79, 11).Value = 24000
$ws.Cells.Item(79, 12).Value = 24000
$ws.Cells.Item(79, 13).Value = 24000
$ws.Cells.Item(79, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(79, 15).Value = "Región del Maule"
$ws.Cells.Item(79, 16).Value = 960
$ws.Cells.Item(79, 17).Value = 25
$ws.Cells.Item(79, 18).Value = "Hortaliza"
